# Insert a new data row at row 187 on the single sheet, shifting the
# existing rows 187:280 down to 188:281 (dimension becomes A1:R281).
# The newly inserted row 187 gets a brand-new data record; every row that
# was previously at N (187 <= N <= 280) now lives at N+1 with its values
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 187 (pushes 187..280 -> 188..281).
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new record's data.
$ws.Cells.Item(187, 1).Value = 6
$ws.Cells.Item(187, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(187, 3).Value = "Metropolitana"
$ws.Cells.Item(187, 4).Value = 44529
$ws.Cells.Item(187, 5).Value = 13
$ws.Cells.Item(187, 6).Value = 100112032
$ws.Cells.Item(187, 7).Value = "Zapallo italiano"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 1010
$ws.Cells.Item(187, 11).Value = 5000
$ws.Cells.Item(187, 12).Value = 6000
$ws.Cells.Item(187, 13).Value = 5446
$ws.Cells.Item(187, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(187, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(187, 16).Value = 109
$ws.Cells.Item(187, 17).Value = 50
$ws.Cells.Item(187, 18).Value = "Hortaliza"
